# Edit script for WEST_VIRGINIA_2024.xlsx
# 1. Rename header columns to snake_case machine-readable names.
# 2. Capitalize connector words ("de", "del", "las", "los", "y") in
#    municipality/state names (Title Case style).
# 3. Change "TOTAL" label in the grand-total row to "Total".
# 4. Remove the trailing metadata rows (155-159) and shrink the used
#    range back down to A1:D153.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row renames -------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case fixes for specific cells --------------------------------
$ws.Range("B8").Value   = "Amatenango De La Frontera"
$ws.Range("B19").Value  = "San Cristóbal De Las Casas"
$ws.Range("A32").Value  = "Ciudad De México"
$ws.Range("A37").Value  = "Coahuila De Zaragoza"
$ws.Range("A44").Value  = "Estado De México"
$ws.Range("B47").Value  = "Naucalpan De Juárez"
$ws.Range("B51").Value  = "Valle De Chalco Solidaridad"
$ws.Range("B61").Value  = "San Miguel De Allende"
$ws.Range("B63").Value  = "Acapulco De Juárez"
$ws.Range("B65").Value  = "Chilapa De Álvarez"
$ws.Range("B66").Value  = "Coyuca De Catalán"
$ws.Range("B70").Value  = "Técpan De Galeana"
$ws.Range("B73").Value  = "Tenango De Doria"
$ws.Range("B74").Value  = "Tepehuacán De Guerrero"
$ws.Range("A80").Value  = "Michoacán De Ocampo"
$ws.Range("B87").Value  = "Tiquicheo De Nicolás Romero"
$ws.Range("B93").Value  = "Putla Villa De Guerrero"
$ws.Range("B102").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B110").Value = "Santo Domingo De Morelos"
$ws.Range("B113").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B114").Value = "Tlacolula De Matamoros"
$ws.Range("B119").Value = "Izúcar De Matamoros"
$ws.Range("B125").Value = "Cadereyta De Montes"
$ws.Range("B126").Value = "Pinal De Amoles"
$ws.Range("B127").Value = "San Juan Del Río"
$ws.Range("A139").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B144").Value = "Ixhuacán De Los Reyes"
$ws.Range("B146").Value = "Poza Rica De Hidalgo"
$ws.Range("B149").Value = "Cañitas De Felipe Pescador"

# --- 3. Grand total label ----------------------------------------------------
$ws.Range("A153").Value = "Total"

# --- 4. Drop the trailing metadata/footer rows (155-159) -------------------
$ws.Range("A155:D159").EntireRow.Delete()
